$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit drops several long narrative paragraphs (Objetivos/Programa
# resumido/Programa/Bibliografia bodies) and the "Docentes responsaveis:"
# label row, moving the four professor lines up so each now sits beside
# the label it ends up under. Net effect: 27 rows collapse to 22.

# --- Rows that need brand-new B/C content where none existed before ---
# (clone the wrap/red-font styling from a row that already has it, then
# overwrite the value, so the new cells pick up styles 2/3 instead of
# inheriting column A's bold style 1)
$ws.Range("B10:C10").Copy($ws.Range("B12:C12"))
$ws.Range("B10:C10").Copy($ws.Range("B18:C18"))
$ws.Range("B10:C10").Copy($ws.Range("B20:C20"))

# --- Row 10: Objetivos: + first professor line (was the long objectives paragraph) ---
$ws.Range("B10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C10").Value = "6495737 - Durval Rodrigues Junior"

# --- Row 11: Objectives: stays label-only ---

# --- Row 12: was "Docentes responsaveis:" label-only -> "Programa resumido:" + professor ---
$ws.Range("A12").Value = "Programa resumido:"
$ws.Range("B12").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("C12").Value = "5983729 - Fernando Vernilli Junior"
$ws.Rows(12).RowHeight = 60

# --- Row 13: was professor-only (Durval) -> "Short syllabus:" label-only ---
$ws.Range("A13").Value = "Short syllabus:"
$ws.Range("B13:C13").Clear()
$ws.Rows(13).RowHeight = 60

# --- Row 14: was professor-only (Fernando) -> "Programa:" + professor (Hugo) ---
$ws.Range("A14").Value = "Programa:"
$ws.Range("B14").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C14").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Rows(14).RowHeight = 120

# --- Row 15: was professor-only (Hugo) -> "Syllabus:" label-only ---
$ws.Range("A15").Value = "Syllabus:"
$ws.Range("B15:C15").Clear()
$ws.Rows(15).RowHeight = 120

# --- Row 16: was professor-only (Maria) -> "Avaliação:" label-only, default height ---
$ws.Range("A16").Value = "Avaliação:"
$ws.Range("B16:C16").Clear()
$ws.Rows(16).UseStandardHeight = $true

# --- Row 17: was "Programa resumido:" + short-syllabus text -> "Método:" + professor (Maria) ---
$ws.Range("A17").Value = "Método:"
$ws.Range("B17").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C17").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Rows(17).RowHeight = 60

# --- Row 18: was "Short syllabus:" label-only -> "Critério:" + exam text ---
$ws.Range("A18").Value = "Critério:"
$ws.Range("B18").Value = "Aplicação de duas provas escritas"
$ws.Range("C18").Value = "Aplicação de duas provas escritas"
$ws.Rows(18).RowHeight = 60

# --- Row 19: was "Programa:" + long program text -> "Norma de recuperação:" + grade formula ---
$ws.Range("A19").Value = "Norma de recuperação:"
$ws.Range("B19").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2"
$ws.Range("C19").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2"
$ws.Rows(19).RowHeight = 60

# --- Row 20: was "Syllabus:" label-only -> "Bibliografia:" + recovery text ---
$ws.Range("A20").Value = "Bibliografia:"
$ws.Range("B20").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C20").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Rows(20).RowHeight = 120

# --- Row 21: was "Avaliação:" label-only -> "Requisitos:" label-only, default height ---
$ws.Range("A21").Value = "Requisitos:"
$ws.Rows(21).UseStandardHeight = $true

# --- Row 22: was "Método:" + exam text -> the prerequisite line (no label) ---
$ws.Range("A22").Clear()
$ws.Range("B22").Value = "LOM3018 -  Introdução à Engenharia de Materiais  (Requisito fraco)`n"
$ws.Range("C22").Value = "LOM3018 -  Introdução à Engenharia de Materiais  (Requisito fraco)`n"
$ws.Rows(22).RowHeight = 30

# Rows 23-27 held the old "Critério/Norma de recuperação/Bibliografia/
# Requisitos/LOM3018" content that has now been folded into rows 18-22
# above, so the trailing rows are no longer needed.
$ws.Rows("23:27").Delete()
